$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook is an odds-data sheet for Uruguay Primera Division matches.
# The update below:
#  - fixes two id-swap mix-ups between rows 114/115 (match ids 7559468/7559469)
#    and a three-way mix-up between rows 117/118/120 (match ids 7013409/7013886/7013702)
#    so each row carries the correct match's data again.
#  - refreshes closing-odds figures for the still-upcoming fixtures in rows 218-222.

# Row 114
$ws.Range("B114").Value = 7559468
$ws.Range("E114").Value = 'Liverpool Montevideo'
$ws.Range("F114").Value = 'CA River Plate'
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 1
$ws.Range("I114").Value = 'H'
$ws.Range("J114").Value = 1.7
$ws.Range("K114").Value = 3
$ws.Range("L114").Value = 5.75
$ws.Range("M114").Value = 1.833
$ws.Range("O114").Value = 4.5
$ws.Range("P114").Value = -0.5
$ws.Range("Q114").Value = 1.925
$ws.Range("R114").Value = 1.925
$ws.Range("S114").Value = 2.25
$ws.Range("T114").Value = 2.025
$ws.Range("U114").Value = 1.825
$ws.Range("V114").Value = 0.833
$ws.Range("W114").Value = -1
$ws.Range("Y114").Value = 0.925
$ws.Range("Z114").Value = -1
$ws.Range("AA114").Value = 1.025
$ws.Range("AB114").Value = -1

# Row 115
$ws.Range("B115").Value = 7559469
$ws.Range("E115").Value = 'Montevideo Wanderers'
$ws.Range("F115").Value = 'Penarol'
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 'D'
$ws.Range("J115").Value = 4.75
$ws.Range("K115").Value = 3.4
$ws.Range("L115").Value = 1.7
$ws.Range("M115").Value = 2.7
$ws.Range("O115").Value = 2.45
$ws.Range("P115").Value = 0
$ws.Range("Q115").Value = 2.05
$ws.Range("R115").Value = 1.8
$ws.Range("S115").Value = 2.5
$ws.Range("T115").Value = 1.975
$ws.Range("U115").Value = 1.875
$ws.Range("V115").Value = -1
$ws.Range("W115").Value = 2.2
$ws.Range("Y115").Value = 0
$ws.Range("Z115").Value = 0
$ws.Range("AA115").Value = -1
$ws.Range("AB115").Value = 0.875

# Row 117
$ws.Range("B117").Value = 7013886
$ws.Range("E117").Value = 'Racing Club de Montevideo'
$ws.Range("F117").Value = 'Cerro'
$ws.Range("G117").Value = 0
$ws.Range("I117").Value = 'A'
$ws.Range("J117").Value = 2.25
$ws.Range("K117").Value = 3.1
$ws.Range("L117").Value = 3.25
$ws.Range("M117").Value = 2.25
$ws.Range("N117").Value = 2.875
$ws.Range("O117").Value = 3.5
$ws.Range("P117").Value = -0.25
$ws.Range("Q117").Value = 1.95
$ws.Range("R117").Value = 1.9
$ws.Range("S117").Value = 2
$ws.Range("T117").Value = 1.925
$ws.Range("U117").Value = 1.925
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 2.5
$ws.Range("Z117").Value = 0.8999999999999999
$ws.Range("AB117").Value = 0.925

# Row 118
$ws.Range("B118").Value = 7013702
$ws.Range("E118").Value = 'Defensor Sporting'
$ws.Range("F118").Value = 'Danubio'
$ws.Range("H118").Value = 2
$ws.Range("J118").Value = 1.8
$ws.Range("K118").Value = 3.6
$ws.Range("L118").Value = 4.2
$ws.Range("M118").Value = 1.8
$ws.Range("N118").Value = 3.6
$ws.Range("O118").Value = 4.2
$ws.Range("P118").Value = -0.75
$ws.Range("Q118").Value = 2.05
$ws.Range("R118").Value = 1.8
$ws.Range("S118").Value = 2.25
$ws.Range("T118").Value = 1.85
$ws.Range("U118").Value = 2
$ws.Range("X118").Value = 3.2
$ws.Range("Z118").Value = 0.8
$ws.Range("AA118").Value = -0.5
$ws.Range("AB118").Value = 0.5

# Row 120
$ws.Range("B120").Value = 7013409
$ws.Range("E120").Value = 'Nacional De Football'
$ws.Range("F120").Value = 'Torque'
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 1
$ws.Range("I120").Value = 'D'
$ws.Range("J120").Value = 1.666
$ws.Range("K120").Value = 3.9
$ws.Range("L120").Value = 4.5
$ws.Range("M120").Value = 1.615
$ws.Range("N120").Value = 4
$ws.Range("O120").Value = 4.75
$ws.Range("Q120").Value = 1.8
$ws.Range("R120").Value = 2.05
$ws.Range("S120").Value = 2.75
$ws.Range("T120").Value = 1.95
$ws.Range("U120").Value = 1.9
$ws.Range("W120").Value = 3
$ws.Range("X120").Value = -1
$ws.Range("Z120").Value = 1.05
$ws.Range("AA120").Value = -1
$ws.Range("AB120").Value = 0.8999999999999999

# Row 218
$ws.Range("M218").Value = 2.875
$ws.Range("Q218").Value = 2.1
$ws.Range("R218").Value = 1.775

# Row 219
$ws.Range("M219").Value = 2.05
$ws.Range("O219").Value = 3.9
$ws.Range("P219").Value = -0.5
$ws.Range("Q219").Value = 2.05
$ws.Range("R219").Value = 1.8
$ws.Range("T219").Value = 1.9
$ws.Range("U219").Value = 1.95

# Row 220
$ws.Range("O220").Value = 2.2
$ws.Range("Q220").Value = 1.975
$ws.Range("R220").Value = 1.875
$ws.Range("T220").Value = 1.85
$ws.Range("U220").Value = 2

# Row 221
$ws.Range("M221").Value = 2.8
$ws.Range("Q221").Value = 1.925
$ws.Range("R221").Value = 1.925
$ws.Range("T221").Value = 1.85
$ws.Range("U221").Value = 2

# Row 222
$ws.Range("M222").Value = 1.85
$ws.Range("O222").Value = 4.2
$ws.Range("Q222").Value = 1.85
$ws.Range("R222").Value = 2
$ws.Range("T222").Value = 2
$ws.Range("U222").Value = 1.85

